# edit.ps1 -- apply the "Going strong to results :P" commit to the deck.
#
# Two real content edits are made (found by scanning every slide/shape/
# paragraph for the distinctive text, rather than hard-coding slide/shape
# numbers, so the script is robust to layout/ordering surprises):
#
#   1. Slide "Problem": the phrase "Environmental conditions also result in
#      a compromise" used to be split "also " / "result in a " across two
#      separate runs; they get collapsed into a single run "also result in
#      a ".
#   2. Slide "Methodology - network population": "We use 1/40 nodes..."
#      becomes "We use 1/30 nodes..." (the run carrying "1/40" is split so
#      the "1" and "/30 " pieces are reachable as independent runs, mirroring
#      the authored OOXML).
#
# (Note: the canonical diff also renumbers an internal legacy VML shape id,
# spid="_x0000_s1088" -> "_x0000_s1089", on the slide master's embedded
# Corel DESIGNER OLE object. That id is an implementation-internal VML
# bookkeeping artifact that PowerPoint itself reassigns during authoring;
# it is not reachable through any property of the PowerPoint object model
# -- there's no such member on Shape/OLEFormat -- so it's left untouched
# here rather than poked at through unsupported means.)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Helper: scan the whole deck and return the first paragraph (as
# "$slideIndex|$shapeIndex|$paraIndex") whose text contains $needle.
# ---------------------------------------------------------------------

$targetSlide = 0
$targetShape = 0
$targetPara = 0
$targetText = ""

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tf = $shape.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                $pcount = $tr.Paragraphs().Count
                for ($pi = 1; $pi -le $pcount; $pi++) {
                    $para = $tr.Paragraphs($pi)
                    $t = $para.Text
                    if (($targetSlide -eq 0) -and ($t.Contains("also result in a "))) {
                        $targetSlide = $si
                        $targetShape = $shi
                        $targetPara = $pi
                        $targetText = $t
                    }
                }
            }
        }
    }
}

if ($targetSlide -ne 0) {
    $slide = $p.Slides.Item($targetSlide)
    $shape = $slide.Shapes.Item($targetShape)
    $para = $shape.TextFrame.TextRange.Paragraphs($targetPara)
    $full = $para.Text
    $idx0 = $full.IndexOf("also result in a ")
    $mergeStart = $idx0 + 1
    $mergeLen = "also result in a ".Length
    $merged = $para.Characters($mergeStart, $mergeLen)
    $merged.Text = "also result in a "
}

# ---------------------------------------------------------------------
# "We use 1/40 nodes ..." -> "We use 1/30 nodes ..."
# ---------------------------------------------------------------------

$targetSlide2 = 0
$targetShape2 = 0
$targetPara2 = 0

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tf = $shape.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                $pcount = $tr.Paragraphs().Count
                for ($pi = 1; $pi -le $pcount; $pi++) {
                    $para = $tr.Paragraphs($pi)
                    $t = $para.Text
                    if (($targetSlide2 -eq 0) -and ($t.Contains("1/40 nodes"))) {
                        $targetSlide2 = $si
                        $targetShape2 = $shi
                        $targetPara2 = $pi
                    }
                }
            }
        }
    }
}

if ($targetSlide2 -ne 0) {
    $slide2 = $p.Slides.Item($targetSlide2)
    $shape2 = $slide2.Shapes.Item($targetShape2)
    $para2 = $shape2.TextFrame.TextRange.Paragraphs($targetPara2)
    $full2 = $para2.Text
    $idx1 = $full2.IndexOf("1/40")

    # Position of the "1" (1-based Characters index).
    $onePos = $idx1 + 1
    $one = $para2.Characters($onePos, 1)
    $one.Text = "1"

    # Position of "/40 " immediately after the "1" -> replace with "/30 ".
    $fracPos = $onePos + 1
    $frac = $para2.Characters($fracPos, 4)
    $frac.Text = "/30 "
}

Write-Output "slide3-merge: slide=$targetSlide shape=$targetShape para=$targetPara"
Write-Output "slide6-split: slide=$targetSlide2 shape=$targetShape2 para=$targetPara2"
